# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
#
# Cell B11 on the "Rules" sheet currently holds the text "R40"; change it
# to hold the text "1" (kept as text, not converted to a number).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new text as a formula result first...
$ws.Range("B11").Formula = "=""1"""
# ...then copy/paste-values-only over itself so the cell becomes a plain
# text value "1" (no formula left behind) without Excel re-interpreting
# the numeric-looking literal as a number (which would also change the
# cell's number format / style).
$ws.Range("B11").Copy()
$ws.Range("B11").PasteSpecial(-4163)
